# Auto-generated edit script: apply "Add data for 2022-09-08" updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("I2").Value = 4976
$ws.Range("I3").Value = 5174
$ws.Range("F4").Value = 1865
$ws.Range("I4").Value = 1185
$ws.Range("I5").Value = 479
$ws.Range("I6").Value = 5652
$ws.Range("F7").Value = 24054
$ws.Range("I7").Value = 17466

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("I2").Value = 137
$ws.Range("I5").Value = 55
$ws.Range("I7").Value = 562
$ws.Range("E8").Value = 1911
$ws.Range("I8").Value = 1062
$ws.Range("I10").Value = 128
$ws.Range("I14").Value = 101
$ws.Range("I15").Value = 203
$ws.Range("I16").Value = 48
$ws.Range("I18").Value = 123
$ws.Range("I19").Value = 479
$ws.Range("I25").Value = 85
$ws.Range("I27").Value = 159
$ws.Range("I29").Value = 1111
$ws.Range("I31").Value = 169
$ws.Range("I33").Value = 797
$ws.Range("I37").Value = 556
$ws.Range("I41").Value = 78
$ws.Range("I42").Value = 593
$ws.Range("I47").Value = 118
$ws.Range("I50").Value = 80
$ws.Range("I51").Value = 187
$ws.Range("I52").Value = 381
$ws.Range("I53").Value = 179
$ws.Range("I55").Value = 195
$ws.Range("I60").Value = 90
$ws.Range("E63").Value = 314
$ws.Range("I63").Value = 67
$ws.Range("I64").Value = 152
$ws.Range("I65").Value = 389
$ws.Range("I67").Value = 683
$ws.Range("I68").Value = 65
$ws.Range("I71").Value = 53
$ws.Range("I76").Value = 265
$ws.Range("F77").Value = 146
$ws.Range("I77").Value = 107
$ws.Range("I78").Value = 249
$ws.Range("I79").Value = 499
$ws.Range("I80").Value = 58
$ws.Range("I85").Value = 791
$ws.Range("I89").Value = 198
$ws.Range("I90").Value = 215
$ws.Range("I95").Value = 286
$ws.Range("I98").Value = 117
$ws.Range("I99").Value = 327
$ws.Range("F101").Value = 24054
$ws.Range("I101").Value = 17466

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("I2").Value = 215
$ws.Range("I7").Value = 791

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("I3").Value = 138
$ws.Range("I7").Value = 381

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("I2").Value = 332
$ws.Range("I3").Value = 298
$ws.Range("E4").Value = 122
$ws.Range("I4").Value = 61
$ws.Range("I5").Value = 29
$ws.Range("I6").Value = 342
$ws.Range("E7").Value = 1911
$ws.Range("I7").Value = 1062

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("I6").Value = 81
$ws.Range("I7").Value = 179

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("I2").Value = 189
$ws.Range("I6").Value = 146
$ws.Range("I7").Value = 562

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("I2").Value = 47
$ws.Range("I6").Value = 73
$ws.Range("I7").Value = 198

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("I2").Value = 36
$ws.Range("I7").Value = 101

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("I2").Value = 175
$ws.Range("I6").Value = 155
$ws.Range("I7").Value = 556

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("I2").Value = 93
$ws.Range("I7").Value = 327

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("I3").Value = 245
$ws.Range("I7").Value = 683

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("I3").Value = 47
$ws.Range("I7").Value = 169

$ws = $wb.Worksheets.Item("New City")
$ws.Range("I3").Value = 116
$ws.Range("I7").Value = 389

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("I2").Value = 98
$ws.Range("I7").Value = 286

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("I3").Value = 297
$ws.Range("I4").Value = 36
$ws.Range("I6").Value = 251
$ws.Range("I7").Value = 797

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("I3").Value = 386
$ws.Range("I6").Value = 301
$ws.Range("I7").Value = 1111

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("I3").Value = 142
$ws.Range("I7").Value = 479

$ws = $wb.Worksheets.Item("River North")
$ws.Range("I6").Value = 119
$ws.Range("I7").Value = 265

$ws = $wb.Worksheets.Item("Hermosa")
$ws.Range("I2").Value = 27
$ws.Range("I7").Value = 78

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("I3").Value = 203
$ws.Range("I6").Value = 167
$ws.Range("I7").Value = 593

$ws = $wb.Worksheets.Item("Avondale")
$ws.Range("I2").Value = 43
$ws.Range("I7").Value = 128

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("I2").Value = 59
$ws.Range("I7").Value = 249

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("I2").Value = 60
$ws.Range("I7").Value = 195

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("I3").Value = 162
$ws.Range("I6").Value = 145
$ws.Range("I7").Value = 499

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("I3").Value = 45
$ws.Range("I7").Value = 152

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("I6").Value = 53
$ws.Range("I7").Value = 123

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("I2").Value = 28
$ws.Range("I7").Value = 85

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("I6").Value = 39
$ws.Range("I7").Value = 118

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("I3").Value = 49
$ws.Range("I7").Value = 203

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("I2").Value = 22
$ws.Range("I7").Value = 117

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("I2").Value = 22
$ws.Range("I7").Value = 80

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("I4").Value = 17
$ws.Range("I7").Value = 137

$ws = $wb.Worksheets.Item("Armour Square")
$ws.Range("I2").Value = 14
$ws.Range("I7").Value = 55

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("I2").Value = 44
$ws.Range("I5").Value = 2
$ws.Range("I7").Value = 159

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("I2").Value = 76
$ws.Range("I3").Value = 48
$ws.Range("I7").Value = 215

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("I6").Value = 74
$ws.Range("I7").Value = 187

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("I3").Value = 20
$ws.Range("I7").Value = 65

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("I6").Value = 27
$ws.Range("I7").Value = 90

$ws = $wb.Worksheets.Item("Oakland")
$ws.Range("I3").Value = 20
$ws.Range("I6").Value = 13
$ws.Range("I7").Value = 53

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("F4").Value = 8
$ws.Range("I4").Value = 5
$ws.Range("F7").Value = 146
$ws.Range("I7").Value = 107

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("I3").Value = 13
$ws.Range("I7").Value = 58

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("I6").Value = 32
$ws.Range("I7").Value = 48
